$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44334
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 9000
$ws.Range("O2").Value = 9000
$ws.Range("P2").Value = 9000
$ws.Range("Q2").Value = '$/caja 15 kilos empedrada'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 600
$ws.Range("T2").Value = 15

$ws.Range("D3").Value = 44314
$ws.Range("L3").Value = 'Especial'
$ws.Range("M3").Value = 210
$ws.Range("N3").Value = 10000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 10000
$ws.Range("Q3").Value = '$/caja 18 kilos granel'
$ws.Range("R3").Value = 'Provincia de Curicó'
$ws.Range("S3").Value = 556
$ws.Range("T3").Value = 18

$ws.Range("D4").Value = 44314
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 140
$ws.Range("N4").Value = 9000
$ws.Range("O4").Value = 9000
$ws.Range("P4").Value = 9000
$ws.Range("Q4").Value = '$/caja 18 kilos granel'
$ws.Range("R4").Value = 'Provincia de Curicó'
$ws.Range("S4").Value = 500
$ws.Range("T4").Value = 18

$ws.Range("D5").Value = 44305
$ws.Range("L5").Value = 'Especial'
$ws.Range("M5").Value = 160
$ws.Range("N5").Value = 10000
$ws.Range("O5").Value = 10000
$ws.Range("P5").Value = 10000
$ws.Range("Q5").Value = '$/caja 18 kilos granel'
$ws.Range("R5").Value = 'Provincia de Curicó'
$ws.Range("S5").Value = 556
$ws.Range("T5").Value = 18

$ws.Range("D6").Value = 44377
$ws.Range("L6").Value = 'Extra (doble especial)'
$ws.Range("M6").Value = 150
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 12000
$ws.Range("Q6").Value = '$/caja 15 kilos empedrada'
$ws.Range("R6").Value = 'Región de O''Higgins'
$ws.Range("S6").Value = 800
$ws.Range("T6").Value = 15

$ws.Range("D7").Value = 44377
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 9000
$ws.Range("O7").Value = 9000
$ws.Range("P7").Value = 9000
$ws.Range("Q7").Value = '$/caja 15 kilos empedrada'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 600
$ws.Range("T7").Value = 15

$ws.Range("D8").Value = 44293
$ws.Range("L8").Value = 'Especial'
$ws.Range("M8").Value = 70
$ws.Range("N8").Value = 10000
$ws.Range("O8").Value = 10000
$ws.Range("P8").Value = 10000
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("R8").Value = 'Provincia de Curicó'
$ws.Range("S8").Value = 556
$ws.Range("T8").Value = 18

$ws.Range("D9").Value = 44293
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 9000
$ws.Range("O9").Value = 9000
$ws.Range("P9").Value = 9000
$ws.Range("Q9").Value = '$/caja 18 kilos granel'
$ws.Range("R9").Value = 'Provincia de Curicó'
$ws.Range("S9").Value = 500
$ws.Range("T9").Value = 18

$ws.Range("D10").Value = 44342
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 300
$ws.Range("N10").Value = 9000
$ws.Range("O10").Value = 9000
$ws.Range("P10").Value = 9000
$ws.Range("Q10").Value = '$/caja 15 kilos empedrada'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 600
$ws.Range("T10").Value = 15

$ws.Range("D11").Value = 44333
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 10
$ws.Range("N11").Value = 130000
$ws.Range("O11").Value = 130000
$ws.Range("P11").Value = 130000
$ws.Range("Q11").Value = '$/bins (450 kilos)'
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value = 289
$ws.Range("T11").Value = 450

$ws.Range("D12").Value = 44424
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 230
$ws.Range("N12").Value = 11000
$ws.Range("O12").Value = 11000
$ws.Range("P12").Value = 11000
$ws.Range("Q12").Value = '$/caja 18 kilos granel'
$ws.Range("R12").Value = 'Región de O''Higgins'
$ws.Range("S12").Value = 611
$ws.Range("T12").Value = 18

$ws.Range("D13").Value = 44294
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 160
$ws.Range("N13").Value = 8000
$ws.Range("O13").Value = 8000
$ws.Range("P13").Value = 8000
$ws.Range("Q13").Value = '$/caja 18 kilos granel'
$ws.Range("R13").Value = 'Provincia de Curicó'
$ws.Range("S13").Value = 444
$ws.Range("T13").Value = 18

$ws.Range("D14").Value = 44312
$ws.Range("L14").Value = 'Especial'
$ws.Range("M14").Value = 230
$ws.Range("N14").Value = 10000
$ws.Range("O14").Value = 10000
$ws.Range("P14").Value = 10000
$ws.Range("Q14").Value = '$/caja 18 kilos granel'
$ws.Range("R14").Value = 'Provincia de Curicó'
$ws.Range("S14").Value = 556
$ws.Range("T14").Value = 18

$ws.Range("D15").Value = 44312
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 9000
$ws.Range("O15").Value = 9000
$ws.Range("P15").Value = 9000
$ws.Range("Q15").Value = '$/caja 18 kilos granel'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 500
$ws.Range("T15").Value = 18

$ws.Range("D16").Value = 44327
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 150
$ws.Range("N16").Value = 10000
$ws.Range("O16").Value = 10000
$ws.Range("P16").Value = 10000
$ws.Range("Q16").Value = '$/caja 18 kilos granel'
$ws.Range("R16").Value = 'Región de O''Higgins'
$ws.Range("S16").Value = 556
$ws.Range("T16").Value = 18

$ws.Range("D17").Value = 44316
$ws.Range("L17").Value = 'Especial'
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 10000
$ws.Range("O17").Value = 10000
$ws.Range("P17").Value = 10000
$ws.Range("Q17").Value = '$/caja 18 kilos granel'
$ws.Range("R17").Value = 'Provincia de Curicó'
$ws.Range("S17").Value = 556
$ws.Range("T17").Value = 18

$ws.Range("D18").Value = 44316
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 180
$ws.Range("N18").Value = 8000
$ws.Range("O18").Value = 8000
$ws.Range("P18").Value = 8000
$ws.Range("Q18").Value = '$/caja 18 kilos granel'
$ws.Range("R18").Value = 'Provincia de Curicó'
$ws.Range("S18").Value = 444
$ws.Range("T18").Value = 18

$ws.Range("D19").Value = 44316
$ws.Range("L19").Value = 'Segunda'
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = 7000
$ws.Range("O19").Value = 7000
$ws.Range("P19").Value = 7000
$ws.Range("Q19").Value = '$/caja 18 kilos granel'
$ws.Range("R19").Value = 'Provincia de Curicó'
$ws.Range("S19").Value = 389
$ws.Range("T19").Value = 18

$ws.Range("D20").Value = 44330
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 8000
$ws.Range("O20").Value = 8000
$ws.Range("P20").Value = 8000
$ws.Range("Q20").Value = '$/caja 15 kilos empedrada'
$ws.Range("R20").Value = 'Región de O''Higgins'
$ws.Range("S20").Value = 533
$ws.Range("T20").Value = 15

$ws.Range("D21").Value = 44330
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 260
$ws.Range("N21").Value = 10000
$ws.Range("O21").Value = 10000
$ws.Range("P21").Value = 10000
$ws.Range("Q21").Value = '$/caja 18 kilos granel'
$ws.Range("R21").Value = 'Región de O''Higgins'
$ws.Range("S21").Value = 556
$ws.Range("T21").Value = 18

$ws.Range("D22").Value = 44280
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 160
$ws.Range("N22").Value = 10000
$ws.Range("O22").Value = 10000
$ws.Range("P22").Value = 10000
$ws.Range("Q22").Value = '$/caja 18 kilos granel'
$ws.Range("R22").Value = 'Provincia de Curicó'
$ws.Range("S22").Value = 556
$ws.Range("T22").Value = 18

$ws.Range("D23").Value = 44323
$ws.Range("L23").Value = 'Especial'
$ws.Range("M23").Value = 300
$ws.Range("N23").Value = 10000
$ws.Range("O23").Value = 10000
$ws.Range("P23").Value = 10000
$ws.Range("Q23").Value = '$/caja 18 kilos granel'
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("S23").Value = 556
$ws.Range("T23").Value = 18

$ws.Range("D24").Value = 44323
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 200
$ws.Range("N24").Value = 8000
$ws.Range("O24").Value = 8000
$ws.Range("P24").Value = 8000
$ws.Range("Q24").Value = '$/caja 18 kilos granel'
$ws.Range("R24").Value = 'Región de O''Higgins'
$ws.Range("S24").Value = 444
$ws.Range("T24").Value = 18

$ws.Range("D25").Value = 44403
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 200
$ws.Range("N25").Value = 10000
$ws.Range("O25").Value = 10000
$ws.Range("P25").Value = 10000
$ws.Range("Q25").Value = '$/caja 18 kilos granel'
$ws.Range("R25").Value = 'Región de O''Higgins'
$ws.Range("S25").Value = 556
$ws.Range("T25").Value = 18

$ws.Range("D26").Value = 44267
$ws.Range("L26").Value = 'Primera'
$ws.Range("M26").Value = 110
$ws.Range("N26").Value = 10000
$ws.Range("O26").Value = 10000
$ws.Range("P26").Value = 10000
$ws.Range("Q26").Value = '$/caja 15 kilos empedrada'
$ws.Range("R26").Value = 'Región de O''Higgins'
$ws.Range("S26").Value = 667
$ws.Range("T26").Value = 15

$ws.Range("D27").Value = 44306
$ws.Range("L27").Value = 'Especial'
$ws.Range("M27").Value = 230
$ws.Range("N27").Value = 10000
$ws.Range("O27").Value = 10000
$ws.Range("P27").Value = 10000
$ws.Range("Q27").Value = '$/caja 18 kilos granel'
$ws.Range("R27").Value = 'Provincia de Curicó'
$ws.Range("S27").Value = 556
$ws.Range("T27").Value = 18

$ws.Range("D28").Value = 44427
$ws.Range("L28").Value = 'Primera'
$ws.Range("M28").Value = 200
$ws.Range("N28").Value = 11000
$ws.Range("O28").Value = 11000
$ws.Range("P28").Value = 11000
$ws.Range("Q28").Value = '$/caja 18 kilos granel'
$ws.Range("R28").Value = 'Región de O''Higgins'
$ws.Range("S28").Value = 611
$ws.Range("T28").Value = 18

$ws.Range("D29").Value = 44291
$ws.Range("L29").Value = 'Primera'
$ws.Range("M29").Value = 300
$ws.Range("N29").Value = 9000
$ws.Range("O29").Value = 9000
$ws.Range("P29").Value = 9000
$ws.Range("Q29").Value = '$/caja 18 kilos granel'
$ws.Range("R29").Value = 'Provincia de Curicó'
$ws.Range("S29").Value = 500
$ws.Range("T29").Value = 18

$ws.Range("D30").Value = 44271
$ws.Range("L30").Value = 'Especial'
$ws.Range("M30").Value = 140
$ws.Range("N30").Value = 10000
$ws.Range("O30").Value = 10000
$ws.Range("P30").Value = 10000
$ws.Range("Q30").Value = '$/caja 15 kilos granel'
$ws.Range("R30").Value = 'Provincia de Curicó'
$ws.Range("S30").Value = 667
$ws.Range("T30").Value = 15

$ws.Range("D31").Value = 44271
$ws.Range("L31").Value = 'Primera'
$ws.Range("M31").Value = 100
$ws.Range("N31").Value = 8000
$ws.Range("O31").Value = 8000
$ws.Range("P31").Value = 8000
$ws.Range("Q31").Value = '$/caja 15 kilos granel'
$ws.Range("R31").Value = 'Provincia de Curicó'
$ws.Range("S31").Value = 533
$ws.Range("T31").Value = 15

$ws.Range("D32").Value = 44286
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 170
$ws.Range("N32").Value = 10000
$ws.Range("O32").Value = 10000
$ws.Range("P32").Value = 10000
$ws.Range("Q32").Value = '$/caja 18 kilos granel'
$ws.Range("R32").Value = 'Provincia de Curicó'
$ws.Range("S32").Value = 556
$ws.Range("T32").Value = 18

$ws.Range("D33").Value = 44350
$ws.Range("L33").Value = 'Especial'
$ws.Range("M33").Value = 150
$ws.Range("N33").Value = 10000
$ws.Range("O33").Value = 10000
$ws.Range("P33").Value = 10000
$ws.Range("Q33").Value = '$/caja 18 kilos granel'
$ws.Range("R33").Value = 'Región de O''Higgins'
$ws.Range("S33").Value = 556
$ws.Range("T33").Value = 18

$ws.Range("D34").Value = 44298
$ws.Range("L34").Value = 'Primera'
$ws.Range("M34").Value = 130
$ws.Range("N34").Value = 9000
$ws.Range("O34").Value = 9000
$ws.Range("P34").Value = 9000
$ws.Range("Q34").Value = '$/caja 18 kilos granel'
$ws.Range("R34").Value = 'Provincia de Curicó'
$ws.Range("S34").Value = 500
$ws.Range("T34").Value = 18

$ws.Range("D35").Value = 44329
$ws.Range("L35").Value = 'Segunda'
$ws.Range("M35").Value = 10
$ws.Range("N35").Value = 130000
$ws.Range("O35").Value = 130000
$ws.Range("P35").Value = 130000
$ws.Range("Q35").Value = '$/bins (450 kilos)'
$ws.Range("R35").Value = 'Región de O''Higgins'
$ws.Range("S35").Value = 289
$ws.Range("T35").Value = 450

$ws.Range("D36").Value = 44302
$ws.Range("L36").Value = 'Primera'
$ws.Range("M36").Value = 310
$ws.Range("N36").Value = 8000
$ws.Range("O36").Value = 9000
$ws.Range("P36").Value = 8452
$ws.Range("Q36").Value = '$/caja 15 kilos empedrada'
$ws.Range("R36").Value = 'Provincia de Curicó'
$ws.Range("S36").Value = 563
$ws.Range("T36").Value = 15

$ws.Range("D37").Value = 44344
$ws.Range("L37").Value = 'Especial'
$ws.Range("M37").Value = 180
$ws.Range("N37").Value = 10000
$ws.Range("O37").Value = 10000
$ws.Range("P37").Value = 10000
$ws.Range("Q37").Value = '$/caja 15 kilos empedrada'
$ws.Range("R37").Value = 'Región de O''Higgins'
$ws.Range("S37").Value = 667
$ws.Range("T37").Value = 15

$ws.Range("D38").Value = 44344
$ws.Range("L38").Value = 'Primera'
$ws.Range("M38").Value = 210
$ws.Range("N38").Value = 9000
$ws.Range("O38").Value = 9000
$ws.Range("P38").Value = 9000
$ws.Range("Q38").Value = '$/caja 15 kilos empedrada'
$ws.Range("R38").Value = 'Región de O''Higgins'
$ws.Range("S38").Value = 600
$ws.Range("T38").Value = 15

$ws.Range("D39").Value = 44357
$ws.Range("L39").Value = 'Primera'
$ws.Range("M39").Value = 230
$ws.Range("N39").Value = 11000
$ws.Range("O39").Value = 11000
$ws.Range("P39").Value = 11000
$ws.Range("Q39").Value = '$/caja 18 kilos granel'
$ws.Range("R39").Value = 'Región de O''Higgins'
$ws.Range("S39").Value = 611
$ws.Range("T39").Value = 18
